# ---------------------------------------------------------------------------
# chore: update Sheets via scheduled runner
#
# Refreshes the cached market-board figures (currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -> columns
# H..N) on each leve-profit worksheet with the latest values pulled by the
# scheduled data-fetch runner. Some leves have no NQ/HQ profit yet (cells
# H:L already both zero) and so have no M/N cell at all; those gain one as
# soon as a profit is computable, and lose it again if it becomes moot.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 17: One for the Road
$ws.Range("H17").Value = 512.5
$ws.Range("J17").Value = 700
$ws.Range("L17").Value = 2100
$ws.Range("N17").Value = -2436

# Row 64: Forged from the Void
$ws.Range("H64").Value = 3842.4285
$ws.Range("I64").Value = 3724.25
$ws.Range("K64").Value = 3724.25
$ws.Range("M64").Value = -3476.25

# Row 67: Dodging the Draft (L)
$ws.Range("H67").Value = 3842.4285
$ws.Range("I67").Value = 3724.25
$ws.Range("K67").Value = 3724.25
$ws.Range("M67").Value = -2866.25

# Row 99: Rumor Has It
$ws.Range("H99").Value = 364
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 599.5
$ws.Range("I137").Value = 599.5
$ws.Range("K137").Value = 1798.5
$ws.Range("M137").Value = 751.5

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2359.1333
$ws.Range("J138").Value = 3916
$ws.Range("L138").Value = 11748
$ws.Range("N138").Value = -22028


# --- ARM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 17503
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 20004
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 20004
$ws.Range("M61").Value = -9788
$ws.Range("N61").Value = -20428

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 2884.8572
$ws.Range("J63").Value = 2719
$ws.Range("L63").Value = 2719
$ws.Range("N63").Value = -4091

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 2884.8572
$ws.Range("J66").Value = 2719
$ws.Range("L66").Value = 13595
$ws.Range("N66").Value = -20459

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1719.4375
$ws.Range("I74").Value = 1700.7333
$ws.Range("K74").Value = 1700.7333
$ws.Range("M74").Value = -826.7333000000001

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1719.4375
$ws.Range("I77").Value = 1700.7333
$ws.Range("K77").Value = 8503.666500000001
$ws.Range("M77").Value = -4135.666500000001

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 3425
$ws.Range("I102").Value = 1850
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 1850
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -228
$ws.Range("N102").Value = -8244

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3311.25
$ws.Range("I132").Value = 3311.25
$ws.Range("K132").Value = 9933.75
$ws.Range("M132").Value = -7403.75

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 17503
$ws.Range("I136").Value = 10000
$ws.Range("J136").Value = 20004
$ws.Range("K136").Value = 30000
$ws.Range("L136").Value = 60012
$ws.Range("M136").Value = -27450
$ws.Range("N136").Value = -65112


# --- BSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 14: Farriers of Fortune
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 1000
$ws.Range("K14").Value = 1000
$ws.Range("M14").Value = -828

# Row 25: Tools of the Trade
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 2315.3333
$ws.Range("I99").Value = 1979.5
$ws.Range("K99").Value = 1979.5
$ws.Range("M99").Value = -481.5

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 1653
$ws.Range("I134").Value = 1653
$ws.Range("K134").Value = 4959
$ws.Range("M134").Value = -2424


# --- CRP ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 6000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 6000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -7248

# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 6000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 30000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -36240

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 4343.2
$ws.Range("J105").Value = 8500
$ws.Range("L105").Value = 8500
$ws.Range("N105").Value = -11994

# Row 109: Playing the Market
$ws.Range("H109").Value = 16999.5
$ws.Range("J109").Value = 16999.5
$ws.Range("L109").Value = 16999.5
$ws.Range("N109").Value = -19079.5

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3499.5
$ws.Range("I134").Value = 3499.5
$ws.Range("K134").Value = 10498.5
$ws.Range("M134").Value = -7963.5


# --- CUL ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 2: Pork Is a Salty Food
$ws.Range("H2").Value = 42.916668
$ws.Range("J2").Value = 109.28571
$ws.Range("L2").Value = 655.71426
$ws.Range("N2").Value = -881.71426

# Row 15: Pretty Enough to Eat
$ws.Range("H15").Value = 145.66667
$ws.Range("I15").Value = 50
$ws.Range("J15").Value = 193.5
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 580.5
$ws.Range("M15").Value = -10
$ws.Range("N15").Value = -860.5

# Row 37: I Love Lamprey
$ws.Range("H37").Value = 200000
$ws.Range("J37").Value = 200000
$ws.Range("L37").Value = 600000
$ws.Range("N37").Value = -600224

# Row 54: Good Eats in Ishgard
$ws.Range("H54").Value = 900
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Row 70: Persona non Gratin
$ws.Range("H70").Value = 30
$ws.Range("I70").Value = 30
$ws.Range("K70").Value = 90
$ws.Range("M70").Value = 225

# Row 73: Recipe for Disaster (L)
$ws.Range("H73").Value = 30
$ws.Range("I73").Value = 30
$ws.Range("K73").Value = 90
$ws.Range("M73").Value = 1002

# Row 114: One Last Meal
$ws.Range("H114").Value = 451.83334
$ws.Range("I114").Value = 624.8333
$ws.Range("J114").Value = 278.83334
$ws.Range("K114").Value = 1874.4999
$ws.Range("L114").Value = 836.5000200000001
$ws.Range("M114").Value = 1379.5001
$ws.Range("N114").Value = -7344.50002

# Row 130: Blast from the Pasta
$ws.Range("H130").Value = 15000
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()


# --- GSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 4396.4
$ws.Range("I97").Value = 4329.3335
$ws.Range("K97").Value = 4329.3335
$ws.Range("M97").Value = -3833.3335

# Row 132: On Board for Lar
$ws.Range("H132").Value = 4062.375
$ws.Range("I132").Value = 3928.5715
$ws.Range("K132").Value = 11785.7145
$ws.Range("M132").Value = -9255.7145


# --- LTW ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 4574.8
$ws.Range("I132").Value = 4527.778
$ws.Range("K132").Value = 13583.334
$ws.Range("M132").Value = -11053.334


# --- WVR ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 96: Skills on Display
$ws.Range("H96").Value = 1750
$ws.Range("I96").Value = 1750
$ws.Range("K96").Value = 1750
$ws.Range("M96").Value = -377

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 2386.6667
$ws.Range("I126").Value = 2386.6667
$ws.Range("K126").Value = 7160.000100000001
$ws.Range("M126").Value = -4690.000100000001

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 1150.6
$ws.Range("I136").Value = 1150.6
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3451.8
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -901.7999999999997
$ws.Range("N136").ClearContents()
